$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Degree symbol, built safely so PowerShell doesn't try to do numeric
# addition with it (e.g. "49.5" + [char]176 gets coerced to a number).
$deg = [char]176

# The at-bat blocks on this sheet got shuffled out of order, and several
# Exit Velo / Launch Angle values that should read "nan" (missing data)
# were showing stale numbers instead. This restores the correct at-bat
# order and fixes the nan placeholders, FB Velo "MPH" suffixes, and the
# pitch-mix abbreviations.

# Block A (rows 10-17)
$ws.Range("J10").Value = 2
$ws.Range("M10").Value = "nan MPH"
$ws.Range("M12").Value = "nan$deg"
$ws.Range("J14").Value = "Roblez"
$ws.Range("M14").Value = "Undefined"
$ws.Range("M15").Value = "Undefined"
$ws.Range("J16").Value = "88-90 MPH"
$ws.Range("J17").Value = "CB,FB,CH"

# Block B (rows 19-26)
$ws.Range("J19").Value = 6
$ws.Range("M19").Value = "92.32 MPH"
$ws.Range("J20").Value = 1
$ws.Range("M21").Value = "49.5$deg"
$ws.Range("M24").Value = "Single"
$ws.Range("J25").Value = "83-85 MPH"
$ws.Range("J26").Value = "SL,CB,FB,CH"

# Block C (rows 28-35)
$ws.Range("J28").Value = 3
$ws.Range("M28").Value = "81.91 MPH"
$ws.Range("J29").Value = 2
$ws.Range("M30").Value = "0.22$deg"
$ws.Range("M32").Value = "Ground Ball"
$ws.Range("M33").Value = "Single"
$ws.Range("J34").Value = "88-90 MPH"
$ws.Range("J35").Value = "CB,FB,CH"

# Block D (partial, rows 43-44)
$ws.Range("J43").Value = "84-84 MPH"
$ws.Range("J44").Value = "SL,FB,CH"

# Block E (rows 46-53)
$ws.Range("J46").Value = 7
$ws.Range("M46").Value = "64.03 MPH"
$ws.Range("M48").Value = "1.21$deg"
$ws.Range("J50").Value = "Plum"
$ws.Range("M50").Value = "Ground Ball"
$ws.Range("M51").Value = "Out"
$ws.Range("J52").Value = "84-86 MPH"
$ws.Range("J53").Value = "SL,FB,CH"

# Block F (rows 61-68)
$ws.Range("J61").Value = 5
$ws.Range("M61").Value = "85.8 MPH"
$ws.Range("J62").Value = 0
$ws.Range("M63").Value = "52.25$deg"
$ws.Range("J65").Value = "Herbst"
$ws.Range("M65").Value = "Fly Ball"
$ws.Range("M66").Value = "Out"
$ws.Range("J67").Value = "83-85 MPH"
$ws.Range("J68").Value = "SL,CB,FB,CH"
